# Update the date line in the title paragraph.
$d = $word.ActiveDocument
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Text = "2024-09-10 Tuesday"

# Update the division problems in the table, addressing each cell directly
# by (row, column) so that duplicate problem text (e.g. "80÷2=" appearing
# twice) is resolved unambiguously and maps to the correct new value.
$tbl = $d.Tables.Item(1)

$tbl.Cell(1,1).Range.Text  = "80÷6="
$tbl.Cell(1,2).Range.Text  = "21÷5="
$tbl.Cell(1,3).Range.Text  = "29÷9="
$tbl.Cell(1,4).Range.Text  = "23÷7="
$tbl.Cell(1,5).Range.Text  = "63÷9="

$tbl.Cell(5,1).Range.Text  = "26÷4="
$tbl.Cell(5,2).Range.Text  = "58÷6="
$tbl.Cell(5,3).Range.Text  = "17÷4="
$tbl.Cell(5,4).Range.Text  = "65÷5="
$tbl.Cell(5,5).Range.Text  = "80÷2="

$tbl.Cell(9,1).Range.Text  = "36÷5="
$tbl.Cell(9,2).Range.Text  = "97÷8="
$tbl.Cell(9,3).Range.Text  = "23÷8="
$tbl.Cell(9,4).Range.Text  = "72÷6="
$tbl.Cell(9,5).Range.Text  = "84÷6="

$tbl.Cell(13,1).Range.Text = "46÷4="
$tbl.Cell(13,2).Range.Text = "60÷8="
$tbl.Cell(13,3).Range.Text = "11÷8="
$tbl.Cell(13,4).Range.Text = "17÷6="
$tbl.Cell(13,5).Range.Text = "19÷8="

$tbl.Cell(17,1).Range.Text = "46÷8="
$tbl.Cell(17,2).Range.Text = "75÷4="
$tbl.Cell(17,3).Range.Text = "36÷9="
$tbl.Cell(17,4).Range.Text = "59÷4="
$tbl.Cell(17,5).Range.Text = "78÷4="
